$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Eintragungsdatum"
$ws.Range("B6").Value = "15.12.2023"

$ws.Range("B5").Select()
